$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new rows right after the header row (old row 2 -> new row 11),
# shifting all existing data rows down by 9.
$ws.Rows("2:10").Insert()

# New data for the freshly-inserted rows 2-10.
$newRows = @(
    @(0.005824529005452867, -0.004876267489825652, 0.009258870057068522),
    @(0.01996676961696429, -0.05454103202494082, -0.001441926153939801),
    @(-0.02569185483247746, -0.01999163068830969, -0.009539442396787617),
    @(0.03522419491999366, 0.005352173823603298, -0.01602810922317025),
    @(-0.01577594932601883, 0.060614168860538, -0.009635333469960495),
    @(-0.03060719080615873, 0.2389583984433218, -0.11525819691028),
    @(-0.02314895105569856, 0.2151985930842024, -0.08631312587233481),
    @(-0.06812567826966884, 0.2275579571723937, -0.07889750547880346),
    @(-0.1647382801355317, 0.1269321128032929, -0.3623354202786154)
)

$r = 2
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# Append one new row of data at the end (row 31).
$ws.Cells.Item(31, 1).Value = 0.02472228522217548
$ws.Cells.Item(31, 2).Value = -0.03295831200341841
$ws.Cells.Item(31, 3).Value = -0.03895686653464332
